# Auto-generated Excel COM-interop script
# Applies updated market-price / profit values to the Leve profit
# tracking tables (one table per game-world job sheet: ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR) as produced by the scheduled price-refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 824.8333
$ws.Cells.Item(12, 9).Value = 812.5
$ws.Cells.Item(12, 10).Value = 849.5
$ws.Cells.Item(12, 11).Value = 812.5
$ws.Cells.Item(12, 12).Value = 849.5
$ws.Cells.Item(12, 13).Value = -642.5
$ws.Cells.Item(12, 14).Value = -1189.5
$ws.Cells.Item(70, 8).Value = 5957.778
$ws.Cells.Item(70, 9).Value = 5124
$ws.Cells.Item(70, 11).Value = 15372
$ws.Cells.Item(70, 13).Value = -15102
$ws.Cells.Item(73, 8).Value = 5957.778
$ws.Cells.Item(73, 9).Value = 5124
$ws.Cells.Item(73, 11).Value = 15372
$ws.Cells.Item(73, 13).Value = -14436
$ws.Cells.Item(86, 8).Value = 4335.9565
$ws.Cells.Item(86, 9).Value = 3042.5
$ws.Cells.Item(86, 11).Value = 3042.5
$ws.Cells.Item(86, 13).Value = -1919.5
$ws.Cells.Item(89, 8).Value = 4335.9565
$ws.Cells.Item(89, 9).Value = 3042.5
$ws.Cells.Item(89, 11).Value = 15212.5
$ws.Cells.Item(89, 13).Value = -9596.5
$ws.Cells.Item(98, 8).Value = 2592.5186
$ws.Cells.Item(98, 9).Value = 2615.3076
$ws.Cells.Item(98, 11).Value = 2615.3076
$ws.Cells.Item(98, 13).Value = -1117.3076
$ws.Cells.Item(100, 8).Value = 12480
$ws.Cells.Item(100, 9).Value = 12480
$ws.Cells.Item(100, 11).Value = 12480
$ws.Cells.Item(100, 13).Value = -11939
$ws.Cells.Item(106, 8).Value = 3399.5
$ws.Cells.Item(106, 9).Value = 3399.5
$ws.Cells.Item(106, 11).Value = 3399.5
$ws.Cells.Item(106, 13).Value = -2768.5
$ws.Cells.Item(116, 8).Value = 5099.636
$ws.Cells.Item(116, 9).Value = 4932.6665
$ws.Cells.Item(116, 10).Value = 5300
$ws.Cells.Item(116, 11).Value = 4932.6665
$ws.Cells.Item(116, 12).Value = 5300
$ws.Cells.Item(116, 13).Value = -1490.6665
$ws.Cells.Item(116, 14).Value = -12184
$ws.Cells.Item(122, 8).Value = 2592.5186
$ws.Cells.Item(122, 9).Value = 2615.3076
$ws.Cells.Item(122, 11).Value = 7845.9228
$ws.Cells.Item(122, 13).Value = -5395.9228
$ws.Cells.Item(138, 8).Value = 3718.2222
$ws.Cells.Item(138, 10).Value = 3374.25
$ws.Cells.Item(138, 12).Value = 10122.75
$ws.Cells.Item(138, 14).Value = -20402.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1528.3572
$ws.Cells.Item(32, 9).Value = 1528.3572
$ws.Cells.Item(32, 11).Value = 1528.3572
$ws.Cells.Item(32, 13).Value = -1241.3572

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3651.1177
$ws.Cells.Item(86, 9).Value = 2998.3333
$ws.Cells.Item(86, 11).Value = 2998.3333
$ws.Cells.Item(86, 13).Value = -1875.3333
$ws.Cells.Item(89, 8).Value = 3651.1177
$ws.Cells.Item(89, 9).Value = 2998.3333
$ws.Cells.Item(89, 11).Value = 14991.6665
$ws.Cells.Item(89, 13).Value = -9375.666499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 469.8
$ws.Cells.Item(22, 9).Value = 516.5
$ws.Cells.Item(22, 10).Value = 399.75
$ws.Cells.Item(22, 11).Value = 516.5
$ws.Cells.Item(22, 12).Value = 399.75
$ws.Cells.Item(22, 13).Value = -166.5
$ws.Cells.Item(22, 14).Value = -1099.75
$ws.Cells.Item(31, 8).Value = 3611.5715
$ws.Cells.Item(31, 9).Value = 2230
$ws.Cells.Item(31, 10).Value = 4916.3887
$ws.Cells.Item(31, 11).Value = 2230
$ws.Cells.Item(31, 12).Value = 4916.3887
$ws.Cells.Item(31, 13).Value = -1935
$ws.Cells.Item(31, 14).Value = -5506.3887
$ws.Cells.Item(34, 8).Value = 3611.5715
$ws.Cells.Item(34, 9).Value = 2230
$ws.Cells.Item(34, 10).Value = 4916.3887
$ws.Cells.Item(34, 11).Value = 2230
$ws.Cells.Item(34, 12).Value = 4916.3887
$ws.Cells.Item(34, 13).Value = -2028
$ws.Cells.Item(34, 14).Value = -5320.3887
$ws.Cells.Item(99, 8).Value = 3814.25
$ws.Cells.Item(99, 10).Value = 4602.8
$ws.Cells.Item(99, 12).Value = 4602.8
$ws.Cells.Item(99, 14).Value = -7598.8
$ws.Cells.Item(122, 8).Value = 1459
$ws.Cells.Item(122, 9).Value = 1580.8572
$ws.Cells.Item(122, 11).Value = 4742.571599999999
$ws.Cells.Item(122, 13).Value = -2292.571599999999
$ws.Cells.Item(126, 8).Value = 3814.25
$ws.Cells.Item(126, 10).Value = 4602.8
$ws.Cells.Item(126, 12).Value = 13808.4
$ws.Cells.Item(126, 14).Value = -18748.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 19998.5
$ws.Cells.Item(80, 10).Value = 19999
$ws.Cells.Item(80, 12).Value = 59997
$ws.Cells.Item(80, 14).Value = -61869
$ws.Cells.Item(83, 8).Value = 19998.5
$ws.Cells.Item(83, 10).Value = 19999
$ws.Cells.Item(83, 12).Value = 179991
$ws.Cells.Item(83, 14).Value = -189351
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 14).Value = ""
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 14).Value = ""
$ws.Cells.Item(107, 8).Value = 2697.9443
$ws.Cells.Item(107, 10).Value = 2628.818
$ws.Cells.Item(107, 12).Value = 7886.454000000001
$ws.Cells.Item(107, 14).Value = -11726.454
$ws.Cells.Item(137, 8).Value = 2558.1875
$ws.Cells.Item(137, 9).Value = 2222.5557
$ws.Cells.Item(137, 11).Value = 6667.6671
$ws.Cells.Item(137, 13).Value = -1567.6671
$ws.Cells.Item(139, 8).Value = 4231.1665
$ws.Cells.Item(139, 9).Value = 3077.4
$ws.Cells.Item(139, 11).Value = 9232.200000000001
$ws.Cells.Item(139, 13).Value = -4092.200000000001
$ws.Cells.Item(140, 8).Value = 5962.5
$ws.Cells.Item(140, 9).Value = 3091.3333
$ws.Cells.Item(140, 11).Value = 9273.999899999999
$ws.Cells.Item(140, 13).Value = -4093.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 90914510
$ws.Cells.Item(80, 9).Value = 200004690
$ws.Cells.Item(80, 10).Value = 6032.1665
$ws.Cells.Item(80, 11).Value = 200004690
$ws.Cells.Item(80, 12).Value = 6032.1665
$ws.Cells.Item(80, 13).Value = -200003692
$ws.Cells.Item(80, 14).Value = -8028.1665
$ws.Cells.Item(83, 8).Value = 90914510
$ws.Cells.Item(83, 9).Value = 200004690
$ws.Cells.Item(83, 10).Value = 6032.1665
$ws.Cells.Item(83, 11).Value = 1000023450
$ws.Cells.Item(83, 12).Value = 30160.8325
$ws.Cells.Item(83, 13).Value = -1000018458
$ws.Cells.Item(83, 14).Value = -40144.8325
$ws.Cells.Item(122, 8).Value = 3574.3157
$ws.Cells.Item(122, 9).Value = 2991.4285
$ws.Cells.Item(122, 11).Value = 8974.2855
$ws.Cells.Item(122, 13).Value = -6524.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1596.8334
$ws.Cells.Item(16, 9).Value = 1014.7273
$ws.Cells.Item(16, 11).Value = 1014.7273
$ws.Cells.Item(16, 13).Value = -844.7273
$ws.Cells.Item(22, 8).Value = 3403.1667
$ws.Cells.Item(22, 10).Value = 2847
$ws.Cells.Item(22, 12).Value = 2847
$ws.Cells.Item(22, 14).Value = -3437
$ws.Cells.Item(26, 8).Value = 8500
$ws.Cells.Item(26, 9).Value = 8500
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 8500
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 13).Value = -8205
$ws.Cells.Item(26, 14).Value = ""
$ws.Cells.Item(27, 8).Value = 3403.1667
$ws.Cells.Item(27, 10).Value = 2847
$ws.Cells.Item(27, 12).Value = 2847
$ws.Cells.Item(27, 14).Value = -3061
$ws.Cells.Item(30, 8).Value = 279.8
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 14).Value = ""
$ws.Cells.Item(46, 8).Value = 3389.7693
$ws.Cells.Item(46, 10).Value = 4926.7
$ws.Cells.Item(46, 12).Value = 4926.7
$ws.Cells.Item(46, 14).Value = -5302.7
$ws.Cells.Item(64, 8).Value = 43000
$ws.Cells.Item(64, 10).Value = 43000
$ws.Cells.Item(64, 12).Value = 43000
$ws.Cells.Item(64, 14).Value = -43450
$ws.Cells.Item(67, 8).Value = 43000
$ws.Cells.Item(67, 10).Value = 43000
$ws.Cells.Item(67, 12).Value = 43000
$ws.Cells.Item(67, 14).Value = -44560
$ws.Cells.Item(100, 8).Value = 2910.7856
$ws.Cells.Item(100, 9).Value = 2035.3
$ws.Cells.Item(100, 11).Value = 2035.3
$ws.Cells.Item(100, 13).Value = -1494.3
$ws.Cells.Item(136, 8).Value = 5835.826
$ws.Cells.Item(136, 9).Value = 1914.6471
$ws.Cells.Item(136, 11).Value = 5743.9413
$ws.Cells.Item(136, 13).Value = -3193.9413
$ws.Cells.Item(137, 8).Value = 58265.8
$ws.Cells.Item(137, 10).Value = 58856.285
$ws.Cells.Item(137, 12).Value = 58856.285
$ws.Cells.Item(137, 14).Value = -69056.285

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1362.55
$ws.Cells.Item(113, 9).Value = 1176.3529
$ws.Cells.Item(113, 11).Value = 3529.0587
$ws.Cells.Item(113, 13).Value = -1359.0587
$ws.Cells.Item(122, 8).Value = 16669259
$ws.Cells.Item(122, 9).Value = 2876.7778
$ws.Cells.Item(122, 11).Value = 8630.3334
$ws.Cells.Item(122, 13).Value = -6180.3334
$ws.Cells.Item(136, 8).Value = 4613.2925
$ws.Cells.Item(136, 9).Value = 4505.077
$ws.Cells.Item(136, 11).Value = 13515.231
$ws.Cells.Item(136, 13).Value = -10965.231
